$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I = I0, J = IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for rows 2-41: column I = I0, column J = IF
$values = @(
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(4, 5),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(5, 6),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(5, 5),
    @(5, 6),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(6, 6),
    @(5, 5),
    @(7, 7),
    @(3, 3),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
